# Apply the refreshed cryptocurrency price/volume snapshot to Sheet1.
# D = Price (text, as scraped -- some values use "." as a thousands separator)
# E = Volume(1h) change, formatted as "  +x.xx%  " (text, fixed leading/trailing spaces)
# A leading "'" forces Excel to keep a numeric-looking Price string as text
# (matches how the source workbook already stores these columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = "43.259.54"
$ws.Range('E2').Value = "  +0.64%  "

# Row 3: Ethereum
$ws.Range('D3').Value = "2.311.29"

# Row 4: TetherUSD
$ws.Range('E4').Value = "  +0.00%  "

# Row 5: BNB
$ws.Range('D5').Value = "'301.75"
$ws.Range('E5').Value = "  +0.07%  "

# Row 6: Solana
$ws.Range('D6').Value = "'98.66"
$ws.Range('E6').Value = "  -0.79%  "

# Row 7: XRP
$ws.Range('E7').Value = "  +3.25%  "

# Row 9: Cardano
$ws.Range('D9').Value = "'0.519"
$ws.Range('E9').Value = "  +0.99%  "

# Row 10: Avalanche
$ws.Range('D10').Value = "'36.49"
$ws.Range('E10').Value = "  +1.21%  "

# Row 11: Dogecoin
$ws.Range('E11').Value = "  +0.35%  "

# Row 12: TRON
$ws.Range('E12').Value = "  +0.54%  "

# Row 13: Chainlink
$ws.Range('D13').Value = "'17.78"
$ws.Range('E13').Value = "  -4.91%  "

# Row 14: Polkadot
$ws.Range('E14').Value = "  -0.74%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = "2.668.24"
$ws.Range('E15').Value = "  +0.44%  "

# Row 16: WrappedEther
$ws.Range('D16').Value = "2.295.49"
$ws.Range('E16').Value = "  -0.15%  "

# Row 17: Polygon
$ws.Range('E17').Value = "  -1.31%  "

# Row 18: WrappedBTC
$ws.Range('D18').Value = "43.122.72"
$ws.Range('E18').Value = "  +0.60%  "

# Row 19: InternetComputer(DFINITY)
$ws.Range('D19').Value = "'13.11"
$ws.Range('E19').Value = "  +4.30%  "

# Row 20: ShibaInu
$ws.Range('D20').Value = "0.0₃0912"
$ws.Range('E20').Value = "  +0.91%  "

# Row 21: Uniswap
$ws.Range('D21').Value = "'6.15"
$ws.Range('E21').Value = "  -0.02%  "

# Row 22: Litecoin
$ws.Range('D22').Value = "'68.42"
$ws.Range('E22').Value = "  +0.82%  "

# Row 23: BitcoinCash
$ws.Range('D23').Value = "'238.61"
$ws.Range('E23').Value = "  +1.33%  "

# Row 24: ImmutableX
$ws.Range('E24').Value = "  -0.04%  "

# Row 25: Dai
$ws.Range('D25').Value = "'1.00"
$ws.Range('E25').Value = "  -0.47%  "

# Row 26: PancakeSwap
$ws.Range('D26').Value = "'2.43"
$ws.Range('E26').Value = "  -0.84%  "

# Row 27: LEO
$ws.Range('E27').Value = "  -0.14%  "

# Row 28: EthereumClassic
$ws.Range('D28').Value = "'25.37"
$ws.Range('E28').Value = "  +1.13%  "

# Row 29: Monero
$ws.Range('D29').Value = "'166.94"
$ws.Range('E29').Value = "  -0.05%  "

# Row 30: Cosmos
$ws.Range('E30').Value = "  +0.51%  "

# Row 31: Toncoin
$ws.Range('E31').Value = "  -11.66%  "

# Row 32: InjectiveProtocol
$ws.Range('D32').Value = "'33.54"
$ws.Range('E32').Value = "  -2.69%  "

# Row 33: Filecoin
$ws.Range('D33').Value = "'5.16"
$ws.Range('E33').Value = "  +2.74%  "

# Row 34: FirstDigitalUSD
$ws.Range('E34').Value = "  -0.04%  "

# Row 35: Celestia
$ws.Range('D35').Value = "'18.20"
$ws.Range('E35').Value = "  +2.44%  "

# Row 36: RenderToken
$ws.Range('D36').Value = "'4.77"
$ws.Range('E36').Value = "  +2.24%  "

# Row 37: WEMIXToken
$ws.Range('E37').Value = "  -0.24%  "

# Row 38: Hedera
$ws.Range('D38').Value = "'0.0694"
$ws.Range('E38').Value = "  +0.42%  "

# Row 39: Kaspa
$ws.Range('E39').Value = "  +1.16%  "

# Row 40: ARBITRUM
$ws.Range('E40').Value = "  +0.41%  "

# Row 41: Stellar
$ws.Range('E41').Value = "  +1.50%  "

# Row 42: LidoDAOToken
$ws.Range('E42').Value = "  -1.54%  "

# Row 43: Maker
$ws.Range('D43').Value = "2.014.22"
$ws.Range('E43').Value = "  +1.49%  "

# Row 44: VeChain
$ws.Range('E44').Value = "  -0.41%  "

# Row 45: ApeXProtocol
$ws.Range('D45').Value = "'2.18"
$ws.Range('E45').Value = "  -6.93%  "

# Row 46: FraxShare
$ws.Range('D46').Value = "'10.33"
$ws.Range('E46').Value = "  +2.01%  "

# Row 47: EnergySwap
$ws.Range('D47').Value = "'17.68"
$ws.Range('E47').Value = "  +0.33%  "

# Row 48: NEARProtocol
$ws.Range('D48').Value = "'2.87"
$ws.Range('E48').Value = "  -0.39%  "

# Row 49: MultiversX
$ws.Range('D49').Value = "'54.59"
$ws.Range('E49').Value = "  -1.45%  "

# Row 50: RocketPoolETH
$ws.Range('D50').Value = "2.538.79"
$ws.Range('E50').Value = "  +0.49%  "

# Row 51: Stacks
$ws.Range('E51').Value = "  +0.46%  "
